$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: "H15"=1656.962, "I15"=1656.962, "K15"=4970.886, "M15"=-4801.886
$ws.Range("H15").Value = 1656.962
$ws.Range("I15").Value = 1656.962
$ws.Range("K15").Value = 4970.886
$ws.Range("M15").Value = -4801.886
# Row 92: "H92"=62500640, "I92"=76923660, "K92"=76923660, "M92"=-76922412
$ws.Range("H92").Value = 62500640
$ws.Range("I92").Value = 76923660
$ws.Range("K92").Value = 76923660
$ws.Range("M92").Value = -76922412
# Row 96: "H96"=16666972, "I96"=16666972, "K96"=50000916, "M96"=-49999543
$ws.Range("H96").Value = 16666972
$ws.Range("I96").Value = 16666972
$ws.Range("K96").Value = 50000916
$ws.Range("M96").Value = -49999543
# Row 98: "H98"=914.2941, "I98"=939.8, "J98"=723, "K98"=939.8, "L98"=723, "M98"=558.2, "N98"=-3719
$ws.Range("H98").Value = 914.2941
$ws.Range("I98").Value = 939.8
$ws.Range("J98").Value = 723
$ws.Range("K98").Value = 939.8
$ws.Range("L98").Value = 723
$ws.Range("M98").Value = 558.2
$ws.Range("N98").Value = -3719
# Row 122: "H122"=914.2941, "I122"=939.8, "J122"=723, "K122"=2819.4, "L122"=2169, "M122"=-369.3999999999996, "N122"=-7069
$ws.Range("H122").Value = 914.2941
$ws.Range("I122").Value = 939.8
$ws.Range("J122").Value = 723
$ws.Range("K122").Value = 2819.4
$ws.Range("L122").Value = 2169
$ws.Range("M122").Value = -369.3999999999996
$ws.Range("N122").Value = -7069
# Row 129: "H129"=180073.67, "J129"=190240.12, "L129"=570720.36, "N129"=-580720.36
$ws.Range("H129").Value = 180073.67
$ws.Range("J129").Value = 190240.12
$ws.Range("L129").Value = 570720.36
$ws.Range("N129").Value = -580720.36
# Row 135: "H135"=12504394, "I135"=422.2, "J135"=100032200, "K135"=3799.8, "L135"=900289800, "M135"=-1264.8, "N135"=-900294870
$ws.Range("H135").Value = 12504394
$ws.Range("I135").Value = 422.2
$ws.Range("J135").Value = 100032200
$ws.Range("K135").Value = 3799.8
$ws.Range("L135").Value = 900289800
$ws.Range("M135").Value = -1264.8
$ws.Range("N135").Value = -900294870
# Row 137: "H137"=25119.596, "I137"=1242.75, "J137"=101525.5, "K137"=3728.25, "L137"=304576.5, "M137"=-1178.25, "N137"=-309676.5
$ws.Range("H137").Value = 25119.596
$ws.Range("I137").Value = 1242.75
$ws.Range("J137").Value = 101525.5
$ws.Range("K137").Value = 3728.25
$ws.Range("L137").Value = 304576.5
$ws.Range("M137").Value = -1178.25
$ws.Range("N137").Value = -309676.5
# Row 138: "H138"=2342.6323, "I138"=2880, "J138"=2249.9827, "K138"=8640, "L138"=6749.9481, "M138"=-3500, "N138"=-17029.9481
$ws.Range("H138").Value = 2342.6323
$ws.Range("I138").Value = 2880
$ws.Range("J138").Value = 2249.9827
$ws.Range("K138").Value = 8640
$ws.Range("L138").Value = 6749.9481
$ws.Range("M138").Value = -3500
$ws.Range("N138").Value = -17029.9481

$ws = $wb.Worksheets.Item("ARM")
# Row 32: "H32"=22483.56, "I32"=25424.918, "J32"=4467.75, "K32"=25424.918, "L32"=4467.75, "M32"=-25137.918, "N32"=-5041.75
$ws.Range("H32").Value = 22483.56
$ws.Range("I32").Value = 25424.918
$ws.Range("J32").Value = 4467.75
$ws.Range("K32").Value = 25424.918
$ws.Range("L32").Value = 4467.75
$ws.Range("M32").Value = -25137.918
$ws.Range("N32").Value = -5041.75
# Row 74: "H74"=37038156, "I74"=52632210, "K74"=52632210, "M74"=-52631336
$ws.Range("H74").Value = 37038156
$ws.Range("I74").Value = 52632210
$ws.Range("K74").Value = 52632210
$ws.Range("M74").Value = -52631336
# Row 77: "H77"=37038156, "I77"=52632210, "K77"=263161050, "M77"=-263156682
$ws.Range("H77").Value = 37038156
$ws.Range("I77").Value = 52632210
$ws.Range("K77").Value = 263161050
$ws.Range("M77").Value = -263156682
# Row 97: "H97"=938.7568, "I97"=948.5517, "J97"=903.25, "K97"=948.5517, "L97"=903.25, "M97"=-452.5517, "N97"=-1895.25
$ws.Range("H97").Value = 938.7568
$ws.Range("I97").Value = 948.5517
$ws.Range("J97").Value = 903.25
$ws.Range("K97").Value = 948.5517
$ws.Range("L97").Value = 903.25
$ws.Range("M97").Value = -452.5517
$ws.Range("N97").Value = -1895.25
# Row 132: "H132"=8557.593999999999, "I132"=1260.1228, "J132"=43220.582, "K132"=3780.3684, "L132"=129661.746, "M132"=-1250.3684, "N132"=-134721.746
$ws.Range("H132").Value = 8557.593999999999
$ws.Range("I132").Value = 1260.1228
$ws.Range("J132").Value = 43220.582
$ws.Range("K132").Value = 3780.3684
$ws.Range("L132").Value = 129661.746
$ws.Range("M132").Value = -1250.3684
$ws.Range("N132").Value = -134721.746
# Row 133: "H133"=0, "J133"=0, "L133"=0
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 135: "H135"=36919.668, "J135"=36919.668, "L135"=36919.668, "N135"=-47059.668
$ws.Range("H135").Value = 36919.668
$ws.Range("J135").Value = 36919.668
$ws.Range("L135").Value = 36919.668
$ws.Range("N135").Value = -47059.668

$ws = $wb.Worksheets.Item("BSM")
# Row 94: "H94"=2183.4546, "I94"=999.75, "K94"=999.75, "M94"=-548.75
$ws.Range("H94").Value = 2183.4546
$ws.Range("I94").Value = 999.75
$ws.Range("K94").Value = 999.75
$ws.Range("M94").Value = -548.75
# Row 99: "H99"=1587.421, "I99"=1610.7858, "J99"=1522, "K99"=1610.7858, "L99"=1522, "M99"=-112.7858000000001, "N99"=-4518
$ws.Range("H99").Value = 1587.421
$ws.Range("I99").Value = 1610.7858
$ws.Range("J99").Value = 1522
$ws.Range("K99").Value = 1610.7858
$ws.Range("L99").Value = 1522
$ws.Range("M99").Value = -112.7858000000001
$ws.Range("N99").Value = -4518
# Row 134: "H134"=79181.86, "I134"=100231.55, "K134"=300694.65, "M134"=-298159.65
$ws.Range("H134").Value = 79181.86
$ws.Range("I134").Value = 100231.55
$ws.Range("K134").Value = 300694.65
$ws.Range("M134").Value = -298159.65

$ws = $wb.Worksheets.Item("CRP")
# Row 22: "H22"=366.66666, "I22"=224.9, "K22"=224.9, "M22"=125.1
$ws.Range("H22").Value = 366.66666
$ws.Range("I22").Value = 224.9
$ws.Range("K22").Value = 224.9
$ws.Range("M22").Value = 125.1
# Row 31: "H31"=10949.757, "I31"=22438.4, "J31"=3116.5908, "K31"=22438.4, "L31"=3116.5908, "M31"=-22143.4, "N31"=-3706.5908
$ws.Range("H31").Value = 10949.757
$ws.Range("I31").Value = 22438.4
$ws.Range("J31").Value = 3116.5908
$ws.Range("K31").Value = 22438.4
$ws.Range("L31").Value = 3116.5908
$ws.Range("M31").Value = -22143.4
$ws.Range("N31").Value = -3706.5908
# Row 34: "H34"=10949.757, "I34"=22438.4, "J34"=3116.5908, "K34"=22438.4, "L34"=3116.5908, "M34"=-22236.4, "N34"=-3520.5908
$ws.Range("H34").Value = 10949.757
$ws.Range("I34").Value = 22438.4
$ws.Range("J34").Value = 3116.5908
$ws.Range("K34").Value = 22438.4
$ws.Range("L34").Value = 3116.5908
$ws.Range("M34").Value = -22236.4
$ws.Range("N34").Value = -3520.5908
# Row 58: "H58"=23981.682, "I58"=1278.9333, "J58"=72630.42999999999, "K58"=1278.9333, "L58"=72630.42999999999, "M58"=-1075.9333, "N58"=-73036.42999999999
$ws.Range("H58").Value = 23981.682
$ws.Range("I58").Value = 1278.9333
$ws.Range("J58").Value = 72630.42999999999
$ws.Range("K58").Value = 1278.9333
$ws.Range("L58").Value = 72630.42999999999
$ws.Range("M58").Value = -1075.9333
$ws.Range("N58").Value = -73036.42999999999
# Row 122: "H122"=1123.0541, "J122"=1109.7727, "L122"=3329.3181, "N122"=-8229.3181
$ws.Range("H122").Value = 1123.0541
$ws.Range("J122").Value = 1109.7727
$ws.Range("L122").Value = 3329.3181
$ws.Range("N122").Value = -8229.3181
# Row 132: "H132"=20007.393, "I132"=23767.783, "J132"=2709.6, "K132"=71303.349, "L132"=8128.799999999999, "M132"=-68773.349, "N132"=-13188.8
$ws.Range("H132").Value = 20007.393
$ws.Range("I132").Value = 23767.783
$ws.Range("J132").Value = 2709.6
$ws.Range("K132").Value = 71303.349
$ws.Range("L132").Value = 8128.799999999999
$ws.Range("M132").Value = -68773.349
$ws.Range("N132").Value = -13188.8
# Row 134: "H134"=740.2759, "I134"=590.34784, "J134"=1315, "K134"=1771.04352, "L134"=3945, "M134"=763.9564799999998, "N134"=-9015
$ws.Range("H134").Value = 740.2759
$ws.Range("I134").Value = 590.34784
$ws.Range("J134").Value = 1315
$ws.Range("K134").Value = 1771.04352
$ws.Range("L134").Value = 3945
$ws.Range("M134").Value = 763.9564799999998
$ws.Range("N134").Value = -9015
# Row 136: "H136"=23981.682, "I136"=1278.9333, "J136"=72630.42999999999, "K136"=3836.7999, "L136"=217891.29, "M136"=-1286.7999, "N136"=-222991.29
$ws.Range("H136").Value = 23981.682
$ws.Range("I136").Value = 1278.9333
$ws.Range("J136").Value = 72630.42999999999
$ws.Range("K136").Value = 3836.7999
$ws.Range("L136").Value = 217891.29
$ws.Range("M136").Value = -1286.7999
$ws.Range("N136").Value = -222991.29

$ws = $wb.Worksheets.Item("CUL")
# Row 131: "H131"=824.77, "J131"=824.77, "L131"=2474.31, "N131"=-12554.31
$ws.Range("H131").Value = 824.77
$ws.Range("J131").Value = 824.77
$ws.Range("L131").Value = 2474.31
$ws.Range("N131").Value = -12554.31
# Row 137: "H137"=27780474, "I137"=1275.8, "J137"=47622760, "K137"=3827.4, "L137"=142868280, "M137"=1272.6, "N137"=-142878480
$ws.Range("H137").Value = 27780474
$ws.Range("I137").Value = 1275.8
$ws.Range("J137").Value = 47622760
$ws.Range("K137").Value = 3827.4
$ws.Range("L137").Value = 142868280
$ws.Range("M137").Value = 1272.6
$ws.Range("N137").Value = -142878480

$ws = $wb.Worksheets.Item("GSM")
# Row 132: "H132"=53292.7, "I132"=56538.156, "J132"=47686.91, "K132"=169614.468, "L132"=143060.73, "M132"=-167084.468, "N132"=-148120.73
$ws.Range("H132").Value = 53292.7
$ws.Range("I132").Value = 56538.156
$ws.Range("J132").Value = 47686.91
$ws.Range("K132").Value = 169614.468
$ws.Range("L132").Value = 143060.73
$ws.Range("M132").Value = -167084.468
$ws.Range("N132").Value = -148120.73

$ws = $wb.Worksheets.Item("LTW")
# Row 82: "H82"=2874.875, "I82"=2499.8333, "J82"=4000, "K82"=2499.8333, "L82"=4000, "M82"=-2138.8333, "N82"=-4722
$ws.Range("H82").Value = 2874.875
$ws.Range("I82").Value = 2499.8333
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 2499.8333
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -2138.8333
$ws.Range("N82").Value = -4722
# Row 85: "H85"=2874.875, "I85"=2499.8333, "J85"=4000, "K85"=2499.8333, "L85"=4000, "M85"=-1251.8333, "N85"=-6496
$ws.Range("H85").Value = 2874.875
$ws.Range("I85").Value = 2499.8333
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 2499.8333
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -1251.8333
$ws.Range("N85").Value = -6496
# Row 93: "H93"=944.3077, "I93"=964.75, "K93"=964.75, "M93"=283.25
$ws.Range("H93").Value = 944.3077
$ws.Range("I93").Value = 964.75
$ws.Range("K93").Value = 964.75
$ws.Range("M93").Value = 283.25
# Row 100: "H100"=3737.375, "I100"=0, "J100"=3737.375, "K100"=0, "L100"=3737.375, "N100"=-4819.375
$ws.Range("H100").Value = 3737.375
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3737.375
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3737.375
$ws.Range("N100").Value = -4819.375
$ws.Range("M100").ClearContents()
# Row 132: "H132"=2753, "I132"=2434.5, "J132"=2912.25, "K132"=7303.5, "L132"=8736.75, "M132"=-4773.5, "N132"=-13796.75
$ws.Range("H132").Value = 2753
$ws.Range("I132").Value = 2434.5
$ws.Range("J132").Value = 2912.25
$ws.Range("K132").Value = 7303.5
$ws.Range("L132").Value = 8736.75
$ws.Range("M132").Value = -4773.5
$ws.Range("N132").Value = -13796.75
# Row 136: "H136"=18068.793, "I136"=21584.791, "K136"=64754.37300000001, "M136"=-62204.37300000001
$ws.Range("H136").Value = 18068.793
$ws.Range("I136").Value = 21584.791
$ws.Range("K136").Value = 64754.37300000001
$ws.Range("M136").Value = -62204.37300000001

$ws = $wb.Worksheets.Item("WVR")
# Row 132: "H132"=978.9091, "I132"=719.5484, "J132"=4999, "K132"=2158.6452, "L132"=14997, "M132"=371.3548000000001, "N132"=-20057
$ws.Range("H132").Value = 978.9091
$ws.Range("I132").Value = 719.5484
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2158.6452
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 371.3548000000001
$ws.Range("N132").Value = -20057
# Row 136: "H136"=23810668, "I136"=27027984, "K136"=81083952, "M136"=-81081402
$ws.Range("H136").Value = 23810668
$ws.Range("I136").Value = 27027984
$ws.Range("K136").Value = 81083952
$ws.Range("M136").Value = -81081402
